# peticiones.xlsx - "eliminado el archivo temporal de excel"
#
# Row 9 used to document a "guardar (cambiar nombre a archivo) archivo de
# datos" endpoint (/archivo/guardar). That temporary/renaming-file endpoint
# is removed and replaced with a new "obtener datos del objeto" (GET all
# data) endpoint. Rows 11-14 referenced the old "/archivo/..." route
# prefix, which is renamed to "/datos/..." to match the new row 9 route.
# Finally the view's scroll position / selection is moved down to the
# bottom of the sheet (A14 top-left, B18 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: replace the "guardar archivo" endpoint with "obtener datos del objeto" ---
$ws.Range("A9").Value = 'obtener datos del objeto'

$ws.Range("B9").Value = '/datos/?version'
$ws.Range("B9").Font.Italic = $true

$ws.Range("C9").Value = 'regresa todos los datos en formato json'
$ws.Range("D9").Value = 'get'
$ws.Range("E9").Value = 'version'

# old F9 (request body) no longer applies to a GET endpoint
$ws.Range("F9").ClearContents()

$ws.Range("G9").Value = '{
 status: ok, 
datos: [{objeto_datos}]
}'

# old H9 was already empty
$ws.Range("H9").ClearContents()

$ws.Rows(9).RowHeight = 46.5

# --- Rows 11-14: "/archivo/..." routes renamed to "/datos/..." ---
$ws.Range("B11").Value = '/datos/?version'
$ws.Range("B12").Value = '/datos/?version/?id'
$ws.Range("B13").Value = '/datos/?version/atributo'
$ws.Range("B14").Value = '/datos/?version/atributo/?nombre'

$ws.Rows(10).RowHeight = 91
$ws.Rows(11).RowHeight = 91
$ws.Rows(13).RowHeight = 57.45
$ws.Rows(14).RowHeight = 57.45
$ws.Rows(15).RowHeight = 57.45
$ws.Rows(16).RowHeight = 57.45
$ws.Rows(17).RowHeight = 57.45
$ws.Rows(18).RowHeight = 57.45

# --- Move the view / selection down to the bottom of the sheet ---
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
[void]$ws.Range("B18").Select()
